# Update automatico via Actualizar 11-10-2020 17-11-10
# Appends 30 new daily rows (44145..44174 / 10-Nov-2020 .. 09-Dec-2020) to the
# UF_IVP_DIARIO sheet, extending the table from row 682 to row 712.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UF_IVP_DIARIO")

$lastRow = 682
$newData = @(
    @(44145, 28895.48, 30112.66),
    @(44146, 28902.2, 30114.99),
    @(44147, 28908.92, 30117.33),
    @(44148, 28915.64, 30119.67),
    @(44149, 28922.37, 30122),
    @(44150, 28929.09, 30124.34),
    @(44151, 28935.82, 30126.68),
    @(44152, 28942.55, 30129.02),
    @(44153, 28949.279999999999, 30131.35),
    @(44154, 28956.01, 30133.69),
    @(44155, 28962.74, 30136.03),
    @(44156, 28969.48, 30138.37),
    @(44157, 28976.22, 30140.71),
    @(44158, 28982.95, 30143.05),
    @(44159, 28989.69, 30145.39),
    @(44160, 28996.44, 30147.72),
    @(44161, 29003.18, 30150.06),
    @(44162, 29009.919999999998, 30152.400000000001),
    @(44163, 29016.67, 30154.74),
    @(44164, 29023.42, 30157.08),
    @(44165, 29030.17, 30159.42),
    @(44166, 29036.92, 30161.759999999998),
    @(44167, 29043.67, 30164.1),
    @(44168, 29050.42, 30166.44),
    @(44169, 29057.18, 30168.79),
    @(44170, 29063.94, 30171.13),
    @(44171, 29070.7, 30173.47),
    @(44172, 29077.46, 30175.81),
    @(44173, 29084.22, 30178.15),
    @(44174, 29090.98, 30180.49)
)

$rowCount = $newData.Count
$firstNewRow = $lastRow + 1
$lastNewRow = $lastRow + $rowCount

# Clone the formatting (number formats / styles) of the last existing data
# row down into the new rows before writing values into them.
$srcRow = $ws.Range("A" + $lastRow + ":C" + $lastRow)
$destRows = $ws.Range("A" + $firstNewRow + ":C" + $lastNewRow)
$srcRow.Copy($destRows)

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $firstNewRow + $i
    $rowVals = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $rowVals[0]
    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
}

# The wider header text no longer needs to wrap onto multiple lines now that
# the columns are widened below, so let the header row's height go back to
# automatic.
$ws.Rows("2").AutoFit()

# Widen columns B and C to fit the (now unwrapped) header text.
$ws.Columns("B").ColumnWidth = 19
$ws.Columns("C").ColumnWidth = 24

# Update the defined name range to cover the newly added rows.
$name = $wb.Names.Item("UF_IVP_DIARIO")
$name.RefersTo = "=UF_IVP_DIARIO!`$A`$1:`$C`$" + $lastNewRow

# Move the selection down to the new last row, matching where the editor's
# cursor ended up after entering the new data.
$ws.Activate()
$ws.Range("B" + $lastNewRow).Select()
